$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# like "0.9990" or "5.120" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.828.00"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "1.885.00"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "239.56"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").Value = "0.4821"
$ws.Range("E7").Value = "  -0.36%  "

$ws.Range("D8").Value = "0.2853"
$ws.Range("E8").Value = "  -1.75%  "

$ws.Range("D9").Value = "0.06551"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").Value = "1.983.00"
$ws.Range("E10").Value = "  +4.51%  "

$ws.Range("D11").Value = "0.07498"
$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("D13").Value = "5.122"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("D14").Value = "88.92"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").Value = "0.6704"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").Value = "30.767.23"
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "2.233.57"
$ws.Range("E17").Value = "  +4.02%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "13.37"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "0.9988"
$ws.Range("E19").Value = "  -0.32%  "

$ws.Range("D20").Value = "0.000007630"
$ws.Range("E20").Value = "  -2.00%  "

$ws.Range("D21").Value = "232.58"
$ws.Range("E21").Value = "  +4.16%  "

$ws.Range("D22").Value = "5.320"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("D23").Value = "0.9993"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").Value = "6.196"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "9.364"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").Value = "168.27"
$ws.Range("E26").Value = "  +2.76%  "

$ws.Range("E27").Value = "  +0.63%  "

$ws.Range("E28").Value = "  +0.30%  "

$ws.Range("D29").Value = "1.418"
$ws.Range("E29").Value = "  -1.95%  "

$ws.Range("D30").Value = "0.09752"
$ws.Range("E30").Value = "  +6.21%  "

$ws.Range("D31").Value = "4.381"
$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").Value = "0.05076"
$ws.Range("E33").Value = "  -0.27%  "

$ws.Range("D34").Value = "1.216"
$ws.Range("E34").Value = "  +4.59%  "

$ws.Range("D35").Value = "0.7559"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").Value = "2.705"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").Value = "0.01876"
$ws.Range("E37").Value = "  -0.97%  "

$ws.Range("D38").Value = "2.635"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").Value = "2.095"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "0.9178"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").Value = "106.64"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("D42").Value = "0.4301"
$ws.Range("E42").Value = "  -1.34%  "

$ws.Range("D43").Value = "5.814"
$ws.Range("E43").Value = "  -3.35%  "

$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  -0.34%  "

$ws.Range("D45").Value = "7.425"
$ws.Range("E45").Value = "  -2.79%  "

$ws.Range("D46").Value = "64.90"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").Value = "0.1290"
$ws.Range("E47").Value = "  -4.08%  "

$ws.Range("D48").Value = "1.486"
$ws.Range("E48").Value = "  -5.90%  "

$ws.Range("D49").Value = "8.936"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "33.97"
$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("D51").Value = "0.3904"
$ws.Range("E51").Value = "  +0.44%  "
